$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.335.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.277.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.592"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.887"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.632.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.323.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.260.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.44%  "
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.74%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0886"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").Value = "  +13.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0330"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +26.89%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.779.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.01%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "60.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "